$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values look numeric but must stay text, like the source data.
$priceCells = @(
    @{Cell="D2"; Value="58.183.16"},
    @{Cell="D3"; Value="2.294.15"},
    @{Cell="D5"; Value="544.30"},
    @{Cell="D6"; Value="131.51"},
    @{Cell="D8"; Value="0.569"},
    @{Cell="D9"; Value="2.288.43"},
    @{Cell="D14"; Value="23.56"},
    @{Cell="D15"; Value="2.699.50"},
    @{Cell="D16"; Value="58.162.11"},
    @{Cell="D17"; Value="0.0000131"},
    @{Cell="D18"; Value="2.270.11"},
    @{Cell="D19"; Value="10.55"},
    @{Cell="D20"; Value="4.27"},
    @{Cell="D21"; Value="312.16"},
    @{Cell="D24"; Value="62.77"},
    @{Cell="D27"; Value="7.97"},
    @{Cell="D30"; Value="170.25"},
    @{Cell="D35"; Value="0.998"},
    @{Cell="D36"; Value="17.69"},
    @{Cell="D40"; Value="37.99"},
    @{Cell="D42"; Value="290.27"},
    @{Cell="D43"; Value="139.74"},
    @{Cell="D45"; Value="0.0949"},
    @{Cell="D48"; Value="18.26"}
)
foreach ($item in $priceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.ClearFormats()
}

# Volume(1h) column (E): plain text percentage strings.
$volumeCells = @(
    @{Cell="E2"; Value="  -3.33%  "},
    @{Cell="E3"; Value="  -5.08%  "},
    @{Cell="E4"; Value="  -0.04%  "},
    @{Cell="E5"; Value="  -2.87%  "},
    @{Cell="E6"; Value="  -3.29%  "},
    @{Cell="E7"; Value="  -0.01%  "},
    @{Cell="E8"; Value="  -3.09%  "},
    @{Cell="E9"; Value="  -5.19%  "},
    @{Cell="E10"; Value="  -4.64%  "},
    @{Cell="E11"; Value="  -2.59%  "},
    @{Cell="E12"; Value="  +0.84%  "},
    @{Cell="E13"; Value="  -5.30%  "},
    @{Cell="E14"; Value="  -5.23%  "},
    @{Cell="E15"; Value="  -5.20%  "},
    @{Cell="E16"; Value="  -3.25%  "},
    @{Cell="E17"; Value="  -4.82%  "},
    @{Cell="E18"; Value="  -5.54%  "},
    @{Cell="E19"; Value="  -5.89%  "},
    @{Cell="E20"; Value="  -5.73%  "},
    @{Cell="E21"; Value="  -4.85%  "},
    @{Cell="E22"; Value="  -5.42%  "},
    @{Cell="E23"; Value="  +0.18%  "},
    @{Cell="E24"; Value="  -3.27%  "},
    @{Cell="E25"; Value="  -4.24%  "},
    @{Cell="E26"; Value="  -0.04%  "},
    @{Cell="E27"; Value="  -7.38%  "},
    @{Cell="E28"; Value="  -6.81%  "},
    @{Cell="E29"; Value="  -3.75%  "},
    @{Cell="E30"; Value="  -0.13%  "},
    @{Cell="E31"; Value="  -7.28%  "},
    @{Cell="E32"; Value="  -1.05%  "},
    @{Cell="E33"; Value="  -7.06%  "},
    @{Cell="E34"; Value="  -6.34%  "},
    @{Cell="E35"; Value="  -0.07%  "},
    @{Cell="E36"; Value="  -4.06%  "},
    @{Cell="E37"; Value="  +0.08%  "},
    @{Cell="E38"; Value="  -8.02%  "},
    @{Cell="E39"; Value="  -6.87%  "},
    @{Cell="E40"; Value="  -1.46%  "},
    @{Cell="E41"; Value="  -6.83%  "},
    @{Cell="E42"; Value="  -10.60%  "},
    @{Cell="E43"; Value="  -5.23%  "},
    @{Cell="E44"; Value="  -5.03%  "},
    @{Cell="E45"; Value="  -2.13%  "},
    @{Cell="E46"; Value="  -3.66%  "},
    @{Cell="E47"; Value="  -4.36%  "},
    @{Cell="E48"; Value="  -7.78%  "},
    @{Cell="E49"; Value="  -3.99%  "},
    @{Cell="E50"; Value="  -0.57%  "},
    @{Cell="E51"; Value="  -5.31%  "}
)
foreach ($item in $volumeCells) {
    $ws.Range($item.Cell).Value = $item.Value
}

